# Append the 02/17/2026 row of profit data (row 85) to Sheet1, matching
# the daily run that appended a new row to data/profit_data.xlsx.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 85

# Column A holds the date as literal text (same style as existing rows,
# e.g. "02/16/2026" in A84) rather than a parsed date serial number.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "02/17/2026"
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).Value = 9410.139999999999
$ws.Cells.Item($row, 3).Value = 0.2423188177825589
$ws.Cells.Item($row, 4).Value = 0.7576811822174411
$ws.Cells.Item($row, 5).Value = -318.78
$ws.Cells.Item($row, 6).Value = -36.02
$ws.Cells.Item($row, 7).Value = -23811.87
$ws.Cells.Item($row, 8).Value = -76.95999999999999
$ws.Cells.Item($row, 9).Value = -1122.6
$ws.Cells.Item($row, 10).Value = -32.99
$ws.Cells.Item($row, 11).Value = -24934.47
$ws.Cells.Item($row, 12).Value = -72.59999999999999
